$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (F, G) before the old "wait" column, which
# shifts the old F ("wait"/2000) data to H.
$ws.Columns("F:G").Insert()

# New F/G columns mirror the existing type/storeValue column pair
# (D = type, E = storeValue) so they pick up the same formatting.
$ws.Range("D1:D3").Copy($ws.Range("F1"))
$ws.Range("E1:E3").Copy($ws.Range("G1"))

# Fill in the new storeValue test step: type id=input1 -> input3.
$ws.Range("G2").ClearContents()
$ws.Range("G3").Value = '{"target":"id=input1","value":"input3"}'
$ws.Range("F3").Value = "input3"

# Row 3 now wraps a long JSON string in G3, matching the row height used
# by rows 1-2.
$ws.Rows("3").RowHeight = 37.5

# Approximate the recalculated (auto-fit) column widths.
$ws.Columns("A").ColumnWidth = 6.65
$ws.Columns("B").ColumnWidth = 64.08
$ws.Columns("C").ColumnWidth = 9.79
$ws.Columns("E").ColumnWidth = 23.36
$ws.Columns("F").ColumnWidth = 9.22
$ws.Columns("G").ColumnWidth = 23.36
$ws.Columns("H").ColumnWidth = 5.5

$ws.Range("C4").Select() | Out-Null
